# Apply a cyclic re-shuffle of several data rows in the "Artfynd" sheet,
# matching the upstream automatic data refresh described in the commit.
#
# Only the cells whose content actually changes are touched, to avoid
# disturbing the representation of any other cell.
#
# Group 1 (rows 6-10): each row's content is replaced by the content that
#   was previously on the next row, wrapping row 10 back to row 6
#   (new(6)=old(7), new(7)=old(8), new(8)=old(9), new(9)=old(10), new(10)=old(6)).
# Group 2 (rows 13-15) and Group 3 (rows 26-28): each row's content is
#   replaced by the content that was previously on the previous row, wrapping
#   the first row of the group back to the last row of the group
#   (new(13)=old(15), new(14)=old(13), new(15)=old(14); and
#    new(26)=old(28), new(27)=old(26), new(28)=old(27)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value2 = 112126260
$ws.Range("B6").Value2 = 77636
$ws.Range("E6").Value2 = 6425
$ws.Range("F6").Value2 = 'Garnlav'
$ws.Range("G6").Value2 = 'Alectoria sarmentosa'
$ws.Range("H6").Value2 = '(Ach.) Ach.'
$ws.Range("Q6").Value2 = 690607
$ws.Range("R6").Value2 = 7125724
$ws.Range("Z6").Value2 = '13:28'
$ws.Range("AB6").Value2 = '13:28'
$ws.Range("A7").Value2 = 112125962
$ws.Range("B7").Value2 = 90844
$ws.Range("E7").Value2 = 5449
$ws.Range("F7").Value2 = 'Svart taggsvamp'
$ws.Range("G7").Value2 = 'Phellodon niger'
$ws.Range("H7").Value2 = '(Fr.:Fr.) P.Karst.'
$ws.Range("Q7").Value2 = 690606
$ws.Range("R7").Value2 = 7125734
$ws.Range("Z7").Value2 = '13:22'
$ws.Range("AB7").Value2 = '13:22'
$ws.Range("A8").Value2 = 112129532
$ws.Range("B8").Value2 = 90800
$ws.Range("D8").Value2 = 'LC'
$ws.Range("E8").Value2 = 4364
$ws.Range("F8").Value2 = 'Dropptaggsvamp'
$ws.Range("G8").Value2 = 'Hydnellum ferrugineum'
$ws.Range("H8").Value2 = '(Fr.:Fr.) P. Karst.'
$ws.Range("P8").Value2 = 'Godmyr (Godmyr), Ly lm'
$ws.Range("Q8").Value2 = 690474
$ws.Range("R8").Value2 = 7126174
$ws.Range("Z8").Value2 = '15:41'
$ws.Range("AB8").Value2 = '15:41'
$ws.Range("A9").Value2 = 112129248
$ws.Range("B9").Value2 = 90844
$ws.Range("D9").Value2 = 'NT'
$ws.Range("E9").Value2 = 5449
$ws.Range("F9").Value2 = 'Svart taggsvamp'
$ws.Range("G9").Value2 = 'Phellodon niger'
$ws.Range("H9").Value2 = '(Fr.:Fr.) P.Karst.'
$ws.Range("Q9").Value2 = 690368
$ws.Range("R9").Value2 = 7126265
$ws.Range("Z9").Value2 = '15:28'
$ws.Range("AB9").Value2 = '15:28'
$ws.Range("A10").Value2 = 112126647
$ws.Range("B10").Value2 = 90434
$ws.Range("E10").Value2 = 4745
$ws.Range("F10").Value2 = 'Tallriska'
$ws.Range("G10").Value2 = 'Lactarius musteus'
$ws.Range("H10").Value2 = 'Fr.'
$ws.Range("P10").Value2 = 'Svarvarmyran (Svarvarmyran), Ly lm'
$ws.Range("Q10").Value2 = 690578
$ws.Range("R10").Value2 = 7125678
$ws.Range("Z10").Value2 = '13:45'
$ws.Range("AB10").Value2 = '13:45'
$ws.Range("A13").Value2 = 112125806
$ws.Range("B13").Value2 = 77636
$ws.Range("E13").Value2 = 6425
$ws.Range("F13").Value2 = 'Garnlav'
$ws.Range("G13").Value2 = 'Alectoria sarmentosa'
$ws.Range("H13").Value2 = '(Ach.) Ach.'
$ws.Range("P13").Value2 = 'Svarvarmyran (Svarvarmyran), Ly lm'
$ws.Range("Q13").Value2 = 690607
$ws.Range("R13").Value2 = 7125748
$ws.Range("S13").Value2 = 2
$ws.Range("Z13").Value2 = '13:16'
$ws.Range("AB13").Value2 = '13:16'
$ws.Range("A14").Value2 = 112129144
$ws.Range("B14").Value2 = 56446
$ws.Range("E14").Value2 = 100049
$ws.Range("F14").Value2 = 'Spillkråka'
$ws.Range("G14").Value2 = 'Dryocopus martius'
$ws.Range("H14").Value2 = '(Linnaeus, 1758)'
$ws.Range("Q14").Value2 = 690342
$ws.Range("R14").Value2 = 7126286
$ws.Range("Z14").Value2 = '15:23'
$ws.Range("AB14").Value2 = '15:23'
$ws.Range("A15").Value2 = 112129344
$ws.Range("B15").Value2 = 90816
$ws.Range("E15").Value2 = 2059
$ws.Range("F15").Value2 = 'Skrovlig taggsvamp'
$ws.Range("G15").Value2 = 'Hydnellum scabrosum'
$ws.Range("H15").Value2 = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("P15").Value2 = 'Godmyr (Godmyr), Ly lm'
$ws.Range("Q15").Value2 = 690448
$ws.Range("R15").Value2 = 7126169
$ws.Range("S15").Value2 = 1
$ws.Range("Z15").Value2 = '15:33'
$ws.Range("AB15").Value2 = '15:33'
$ws.Range("A26").Value2 = 112129437
$ws.Range("B26").Value2 = 90792
$ws.Range("D26").Value2 = 'NT'
$ws.Range("E26").Value2 = 4361
$ws.Range("F26").Value2 = 'Orange taggsvamp'
$ws.Range("G26").Value2 = 'Hydnellum aurantiacum'
$ws.Range("H26").Value2 = '(Batsch:Fr.) P.Karst.'
$ws.Range("P26").Value2 = 'Godmyr (Godmyr), Ly lm'
$ws.Range("Q26").Value2 = 690460
$ws.Range("R26").Value2 = 7126186
$ws.Range("Z26").Value2 = '15:36'
$ws.Range("AB26").Value2 = '15:36'
$ws.Range("A27").Value2 = 112127162
$ws.Range("B27").Value2 = 90800
$ws.Range("D27").Value2 = 'LC'
$ws.Range("E27").Value2 = 4364
$ws.Range("F27").Value2 = 'Dropptaggsvamp'
$ws.Range("G27").Value2 = 'Hydnellum ferrugineum'
$ws.Range("H27").Value2 = '(Fr.:Fr.) P. Karst.'
$ws.Range("Q27").Value2 = 690562
$ws.Range("R27").Value2 = 7125622
$ws.Range("Z27").Value2 = '14:06'
$ws.Range("AB27").Value2 = '14:06'
$ws.Range("A28").Value2 = 112127443
$ws.Range("B28").Value2 = 90843
$ws.Range("E28").Value2 = 5448
$ws.Range("F28").Value2 = 'Svartvit taggsvamp'
$ws.Range("G28").Value2 = 'Phellodon connatus'
$ws.Range("H28").Value2 = '(Schultz) nom.prov'
$ws.Range("P28").Value2 = 'Svarvarmyran (Svarvarmyran), Ly lm'
$ws.Range("Q28").Value2 = 690559
$ws.Range("R28").Value2 = 7125584
$ws.Range("Z28").Value2 = '14:21'
$ws.Range("AB28").Value2 = '14:21'
